$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply an AutoFilter on column D ("Component", colId=3) of Table1's
#     range, keeping only rows whose value is "Type". This both writes the
#     <filterColumn><filters><filter val="Type"/></filters></filterColumn>
#     into the table definition AND hides the non-matching data rows
#     (rows 2-51) exactly like Excel does when a filter is applied.
$ws.Range("A1:O93").AutoFilter(4, @("Type"), 7)

# --- Update the "Spec" column (I) cells that referenced the old shared
#     string "NDR 4.0" to use the new "NDR-4.0" value for the rows that
#     now use the new niem-model type defs spec naming.
$specRows = @(52,53,55,63,64,65,66,69,72,73,74,77,78,79,83,84,86)
foreach ($r in $specRows) {
    $ws.Cells.Item($r, 9).Value = "NDR-4.0"
}

# --- Move the value in M57 (Exceptions) over to N57 (Exception IDs).
$m57 = $ws.Range("M57").Value()
$ws.Range("N57").Value = $m57
$ws.Range("M57").Clear()

# --- Update the frozen-pane top-left cell and the active selection in the
#     bottom-right pane to reflect the new scroll position / selected cell.
$ws.Range("I52").Select()
